$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Menu" sheet - added right after "CarryOut" (the last sheet at the time)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsMenu = $wb.Worksheets.Add($null, $lastSheet)
$wsMenu.Name = "Menu"

$wsMenu.Range("A1").Value = "Category"
$wsMenu.Range("A2").Value = "Appetizers"
$wsMenu.Range("B1").Value = "Item"
$wsMenu.Range("B2").Value = "Southwestern Eggrolls"

$wsMenu.Columns("A:B").AutoFit() | Out-Null
[void]$wsMenu.Range("D13").Select()

# ---------------------------------------------------------------------------
# 2) "DeliveryASAP" sheet - added right after "Menu"
# ---------------------------------------------------------------------------
$wsDelivery = $wb.Worksheets.Add($null, $wsMenu)
$wsDelivery.Name = "DeliveryASAP"

$wsDelivery.Range("A1").Value = "First Name"
$wsDelivery.Range("B1").Value = "Last Name"
$wsDelivery.Range("D1").Value = "Email"

$wsDelivery.Range("A2").Value = "Test1"
$wsDelivery.Range("B2").Value = "Data1"
$wsDelivery.Range("C1").Value = "Contact Number"

$wsDelivery.Range("D2").Value = "abcd@qmail.com"
$wsDelivery.Hyperlinks.Add($wsDelivery.Range("D2"), "mailto:abcd@qmail.com") | Out-Null
$wsDelivery.Range("D2").Style = "Hyperlink"

[void]$wsDelivery.Range("E10").Select()

# ---------------------------------------------------------------------------
# 3) "GuestUserCurbSide" sheet - added right after "DeliveryASAP"
# ---------------------------------------------------------------------------
$wsGuest = $wb.Worksheets.Add($null, $wsDelivery)
$wsGuest.Name = "GuestUserCurbSide"

$wsGuest.Range("A1").Value = "Vehicle Make"
$wsGuest.Range("A2").Value = "Toyota"
$wsGuest.Range("B1").Value = "Vehicle Model"
$wsGuest.Range("B2").Value = "E1346"
$wsGuest.Range("C1").Value = "Vehicle Color"
$wsGuest.Range("C2").Value = "Blue"

$wsGuest.Columns("A:C").AutoFit() | Out-Null
[void]$wsGuest.Range("H13").Select()

# ---------------------------------------------------------------------------
# 4) back to "DeliveryASAP" for the last cell (matches shared-string order)
# ---------------------------------------------------------------------------
$wsDelivery.Range("C2").Value = "(512)242-3434"
$wsDelivery.Columns("A:D").AutoFit() | Out-Null

$wsDelivery.Activate()
[void]$wsDelivery.Range("E10").Select()
